$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.378.65"
$ws.Range("E2").Value = "  -4.35%  "
$ws.Range("D3").Value = "1.567.23"
$ws.Range("E3").Value = "  -4.52%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "289.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3690"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.50%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.43"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.77%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3388"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.168"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07593"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.11"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.048"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.880"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.58%  "
$ws.Range("D16").Value = "1.573.57"
$ws.Range("E16").Value = "  -4.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001132"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "89.08"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -7.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06756"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.86%  "
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.227"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.90%  "
$ws.Range("B22").Value = "BitDAO"
$ws.Range("C22").Value = "https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.5329"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.64%  "
$ws.Range("B23").Value = "Avalanche"
$ws.Range("C23").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "16.51"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.70%  "
$ws.Range("D25").Value = "22.397.15"
$ws.Range("E25").Value = "  -4.38%  "
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.000"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.16%  "
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.375"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.30%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.88"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "145.66"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.968"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "125.24"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.55%  "
$ws.Range("D32").Value = "1.748.10"
$ws.Range("E32").Value = "  -4.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.051"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.264"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.991"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.78%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.28"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -9.86%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08478"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02544"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2333"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.71%  "
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.531"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.72%  "
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.06519"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.80"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -8.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.245"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.83%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6378"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.68%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5982"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.49%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.767"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.131"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.92%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.264"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "123.37"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.90%  "
